# Product Backlog: rename the original sheet to "Sprint 1" and add a new,
# empty "Sprint 2" sheet right after it for next sprint's planning, then
# leave the view on Sprint 1 scrolled/zoomed the way the author left it
# (zoomed out to 70% with the selection on D27).

$wb = $excel.ActiveWorkbook

# Rename the existing (only) worksheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sprint 1"

# Add the new sheet for the next sprint, placed right after Sprint 1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sprint 2"

# Keep Sprint 1 as the active/selected tab, adjust zoom and selection.
$ws1.Activate()
$excel.ActiveWindow.Zoom = 70
$ws1.Range("D27").Select() | Out-Null
